$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("5010754-58.2017.8.21.0001", "0196807-38.2017.8.21.0001", "Digitalizado"),
    @("5008281-41.2013.8.21.0001", "0087776-25.2013.8.21.0001", "Digitalizado"),
    @("5001221-76.2018.8.21.0054", "9000149-83.2018.8.21.0054", "Migrado"),
    @("5006432-89.2018.8.21.0023", "9003986-45.2018.8.21.0023", "Migrado"),
    @("5006429-37.2018.8.21.0023", "9003787-23.2018.8.21.0023", "Migrado"),
    @("5009614-96.2011.8.21.0001", "0420415-91.2011.8.21.0001", "Digitalizado"),
    @("5033806-25.2013.8.21.0001", "0007850-16.2013.8.21.3001", "Digitalizado")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
